# Update the "取得日時" (retrieved timestamp) column on the active sheet
# (ランサーズ) from 2025-12-16 06:30:30 to 2025-12-16 06:39:28 for every
# data row currently present (rows 2-11), matching the commit
# "Append: 2025-12-16 06:39 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2025-12-16 06:30:30"
$newTimestamp = "2025-12-16 06:39:28"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
